$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.856.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.374.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.798.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.801.74'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.369.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.48'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '320.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  -3.05%  '
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  -2.15%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('E29').Value = '  -1.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('E32').Value = '  +11.43%  '
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '320.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '145.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0966'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.65'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.942'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
